$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range runs) ---
# si#6 runs concatenate to "Volume 30   Number  9" -> "...10"
$ws.Range("A8").Value = "Volume 30   Number  10"
# si#9 runs concatenate to "Report Covering the Week  2/27/2023  Through  3/5/2023" -> new dates
$ws.Range("C9").Value = "Report Covering the Week  3/6/2023  Through  3/12/2023"

# --- Weekly crime-stat figures refresh (rows 14-29) ---
# Plain numeric value updates (style/type already correct, no donor-paste needed).
# Row 14
$ws.Range("M14").Value = -50
$ws.Range("N14").Value = -91.666666666666
# Row 15
$ws.Range("C15").Value = 2
$ws.Range("F15").Value = 4
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 300
$ws.Range("I15").Value = 6
$ws.Range("K15").Value = 50
$ws.Range("L15").Value = 50
$ws.Range("M15").Value = 500
$ws.Range("N15").Value = -25
# Row 16
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = 40
$ws.Range("G16").Value = 25
$ws.Range("H16").Value = -16
$ws.Range("I16").Value = 60
$ws.Range("J16").Value = 70
$ws.Range("K16").Value = -14.285714285714
$ws.Range("L16").Value = -6.25
$ws.Range("M16").Value = -15.492957746478
$ws.Range("N16").Value = -86.899563318777
# Row 17
$ws.Range("C17").Value = 11
$ws.Range("D17").Value = 16
$ws.Range("E17").Value = -31.25
$ws.Range("F17").Value = 40
$ws.Range("G17").Value = 41
$ws.Range("H17").Value = -2.439024390243
$ws.Range("I17").Value = 123
$ws.Range("J17").Value = 124
$ws.Range("K17").Value = -0.806451612903
$ws.Range("L17").Value = 33.695652173913
$ws.Range("M17").Value = 73.239436619718
$ws.Range("N17").Value = -30.113636363636
# Row 18
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -25
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = 8.333333333333
$ws.Range("I18").Value = 33
$ws.Range("J18").Value = 43
$ws.Range("K18").Value = -23.255813953488
$ws.Range("L18").Value = -25
$ws.Range("M18").Value = -21.428571428571
$ws.Range("N18").Value = -92.009685230024
# Row 19
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = -50
$ws.Range("F19").Value = 30
$ws.Range("G19").Value = 48
$ws.Range("H19").Value = -37.5
$ws.Range("I19").Value = 95
$ws.Range("J19").Value = 121
$ws.Range("K19").Value = -21.487603305785
$ws.Range("L19").Value = 1.063829787234
$ws.Range("M19").Value = 102.127659574468
$ws.Range("N19").Value = -46.022727272727
# Row 20
$ws.Range("C20").Value = 9
$ws.Range("D20").Value = 9
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 19
$ws.Range("G20").Value = 19
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 55
$ws.Range("J20").Value = 57
$ws.Range("K20").Value = -3.508771929824
$ws.Range("L20").Value = 223.529411764706
$ws.Range("M20").Value = 139.130434782609
$ws.Range("N20").Value = -74.537037037037
# Row 21
$ws.Range("C21").Value = 40
$ws.Range("D21").Value = 50
$ws.Range("E21").Value = -20
$ws.Range("F21").Value = 128
$ws.Range("G21").Value = 146
$ws.Range("H21").Value = -12.328767123287
$ws.Range("I21").Value = 373
$ws.Range("J21").Value = 420
$ws.Range("K21").Value = -11.190476190476
$ws.Range("L21").Value = 18.037974683544
$ws.Range("M21").Value = 45.136186770428
$ws.Range("N21").Value = -74.434544208361
# Row 22
$ws.Range("D22").Value = 1
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -33.333333333333
$ws.Range("J22").Value = 7
$ws.Range("K22").Value = -28.571428571428
# Row 23
$ws.Range("M23").Value = 20
# Row 24
$ws.Range("C24").Value = 14
$ws.Range("D24").Value = 23
$ws.Range("E24").Value = -39.130434782608
$ws.Range("F24").Value = 59
$ws.Range("G24").Value = 108
$ws.Range("H24").Value = -45.370370370370
$ws.Range("I24").Value = 164
$ws.Range("J24").Value = 241
$ws.Range("K24").Value = -31.950207468879
$ws.Range("L24").Value = 3.797468354430
$ws.Range("M24").Value = -13.227513227513
# Row 25
$ws.Range("C25").Value = 19
$ws.Range("D25").Value = 17
$ws.Range("E25").Value = 11.764705882352
$ws.Range("F25").Value = 87
$ws.Range("G25").Value = 84
$ws.Range("H25").Value = 3.571428571428
$ws.Range("I25").Value = 235
$ws.Range("J25").Value = 184
$ws.Range("K25").Value = 27.717391304347
$ws.Range("L25").Value = 52.597402597402
$ws.Range("M25").Value = 12.980769230769
# Row 26
$ws.Range("C26").Value = 2
$ws.Range("F26").Value = 4
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 13
$ws.Range("K26").Value = 44.444444444444
$ws.Range("L26").Value = 85.714285714285
# Row 27
$ws.Range("C27").Value = 2
$ws.Range("F27").Value = 7
$ws.Range("G27").Value = 9
$ws.Range("H27").Value = -22.222222222222
$ws.Range("I27").Value = 17
$ws.Range("K27").Value = -10.526315789473
$ws.Range("L27").Value = 30.769230769230
# Row 28
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 2
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 7
$ws.Range("K28").Value = -14.285714285714
# Row 29
$ws.Range("D29").Value = 2
$ws.Range("E29").Value = -100
$ws.Range("F29").Value = 2
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 7
$ws.Range("K29").Value = -28.571428571428

# --- Type-swap cells: numeric <-> text-placeholder ("0" / "***.*") ---
# These need NumberFormat switched to "@"/"General" + value set, THEN a
# Copy/PasteSpecial(xlPasteFormats) from a donor cell that already carries
# the exact target style index, so the saved style id matches exactly
# (not a newly minted "General"/"@" style).
$xlPasteFormats = -4122

# Donor cells (unaffected by this edit, row 14):
#   C14 = style 14, text "0"      E14 = style 14, text "***.*"
#   F14 = style 15, plain number  K14 = style 16, decimal-format number

# numeric -> text "0"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D15").PasteSpecial($xlPasteFormats)

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D26").PasteSpecial($xlPasteFormats)

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D27").PasteSpecial($xlPasteFormats)

$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C28").PasteSpecial($xlPasteFormats)

$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C29").PasteSpecial($xlPasteFormats)

# numeric -> text "***.*"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E15").PasteSpecial($xlPasteFormats)

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E26").PasteSpecial($xlPasteFormats)

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E27").PasteSpecial($xlPasteFormats)

# text -> plain number (style 15)
$ws.Range("D20").Value = 9
$ws.Range("F14").Copy()
$ws.Range("D20").PasteSpecial($xlPasteFormats)

$ws.Range("D28").Value = 2
$ws.Range("F14").Copy()
$ws.Range("D28").PasteSpecial($xlPasteFormats)

$ws.Range("D29").Value = 2
$ws.Range("F14").Copy()
$ws.Range("D29").PasteSpecial($xlPasteFormats)

# text -> decimal-format number (style 16)
$ws.Range("E20").Value = 0
$ws.Range("K14").Copy()
$ws.Range("E20").PasteSpecial($xlPasteFormats)

$ws.Range("E28").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E28").PasteSpecial($xlPasteFormats)

$ws.Range("E29").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E29").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false